$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (F:G) for the "icp" group, shifting the existing
# mlag/vlag/mlos/vlos/Comment columns two places to the right.
$ws.Range("F1:G1").EntireColumn.Insert()

# Rename the old "pic" header labels to "icp" (same position, new text).
$ws.Range("D1").Value = "micp"
$ws.Range("E1").Value = "vicp"

# Header labels for the brand-new "adp" columns.
$ws.Range("F1").Value = "madp"
$ws.Range("G1").Value = "vadp"

# Fill the new "adp" data columns (same value, 1 / 0, for every data row).
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 6).Value = 1
    $ws.Cells.Item($row, 7).Value = 0
}

# Match the author's final cursor position.
$ws.Range("F12").Select() | Out-Null
